# Update the "Förändrad" (Changed) date column C for rows 2-7 from 2023-10-22 (45221) to 2023-10-25 (45224)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
